$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Overall table preferred width (w:tblW): 11334 dxa -> 9169 dxa
$t.PreferredWidth = 458.45

# Table indent from left (w:tblInd): -1466 dxa -> -336 dxa
$t.Rows.LeftIndent = -16.8

# Column widths (w:gridCol / w:tcW), all rows share the same grid
$t.Columns.Item(1).Width = 95.3   # 2356 -> 1906 dxa
$t.Columns.Item(2).Width = 69.05  # 1707 -> 1381 dxa
$t.Columns.Item(3).Width = 69.05  # 1707 -> 1381 dxa
$t.Columns.Item(4).Width = 69.05  # 1707 -> 1381 dxa
$t.Columns.Item(5).Width = 60.1   # 1486 -> 1202 dxa
$t.Columns.Item(6).Width = 95.9   # 2371 -> 1918 dxa

# Row heights (w:trHeight)
$t.Rows.Item(1).Height = 53.95  # 1176 -> 1079 dxa
$t.Rows.Item(2).Height = 27.7   # 604  -> 554  dxa
$t.Rows.Item(3).Height = 25.95  # 566  -> 519  dxa
